$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This shared string is used by Overview!E2, Overview!F2, zh-cn!C2 and de-de!C2,
# so all four cells need to show the new text.
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- zh-cn row 2: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$mdName = $zhcn.Range("A2").Text
$zhHandoffFile = $zhcn.Range("G2").Text

$zhcn.Range("I2").Value = $mdName
$zhcn.Range("J2").Value = $zhHandoffFile
$zhcn.Range("K2").Value = "2016-08-20 11:03:00"

$zhLink = $zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d2583ce0a8c799708a297aee1911b313ecbf879/e2e/f5a77df7-f7ff-4e94-90a5-abe61fdc42c9.md", [Type]::Missing, [Type]::Missing, $mdName)
$zhFont = $zhcn.Range("I2").Font
$zhFont.Underline = 2
$zhFont.Color = 15570276

# --- de-de row 2: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$mdNameDe = $dede.Range("A2").Text
$deHandoffFile = $dede.Range("G2").Text

$dede.Range("I2").Value = $mdNameDe
$dede.Range("J2").Value = $deHandoffFile
$dede.Range("K2").Value = "2016-08-20 11:03:11"

$deLink = $dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d2583ce0a8c799708a297aee1911b313ecbf879/e2e/f5a77df7-f7ff-4e94-90a5-abe61fdc42c9.md", [Type]::Missing, [Type]::Missing, $mdNameDe)
$deFont = $dede.Range("I2").Font
$deFont.Underline = 2
$deFont.Color = 15570276
